# Lattice multiplication worksheet update
# - regenerates the 15 exercises in the 5x3 table
# - row 4 of the original table ("94 x 66", "46 x 34", "50 x 15") is removed
# - a new row ("25 x 18", "70 x 31", "56 x 54") is appended at the end

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$br = [char]11

function Set-CellContent($table, $row, $col, $top, $digits, $l1, $l2) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $top + $br + $digits + $br + "  ----" + $br + $l1 + $br + $l2
}

# Remove the old fourth row ("94 x 66" / "46 x 34" / "50 x 15") entirely.
$t.Rows.Item(4).Delete()

# Row 1
Set-CellContent $t 1 1 "62 x 69" "  6    9" "6|    |" "2|    |"
Set-CellContent $t 1 2 "38 x 98" "  9    8" "3|    |" "8|    |"
Set-CellContent $t 1 3 "30 x 33" "  3    3" "3|    |" "0|    |"

# Row 2
Set-CellContent $t 2 1 "70 x 42" "  4    2" "7|    |" "0|    |"
Set-CellContent $t 2 2 "19 x 59" "  5    9" "1|    |" "9|    |"
Set-CellContent $t 2 3 "93 x 77" "  7    7" "9|    |" "3|    |"

# Row 3
Set-CellContent $t 3 1 "63 x 61" "  6    1" "6|    |" "3|    |"
Set-CellContent $t 3 2 "13 x 89" "  8    9" "1|    |" "3|    |"
Set-CellContent $t 3 3 "70 x 11" "  1    1" "7|    |" "0|    |"

# Row 4 (originally row 5 before the deletion above)
Set-CellContent $t 4 1 "81 x 36" "  3    6" "8|    |" "1|    |"
Set-CellContent $t 4 2 "28 x 38" "  3    8" "2|    |" "8|    |"
Set-CellContent $t 4 3 "76 x 32" "  3    2" "7|    |" "6|    |"

# New row 5, appended at the end of the table.
$t.Rows.Add() | Out-Null
$newRowIdx = $t.Rows.Count
Set-CellContent $t $newRowIdx 1 "25 x 18" "  1    8" "2|    |" "5|    |"
Set-CellContent $t $newRowIdx 2 "70 x 31" "  3    1" "7|    |" "0|    |"
Set-CellContent $t $newRowIdx 3 "56 x 54" "  5    4" "5|    |" "6|    |"
